# Workbook was edited on 06/03: the "--LIMIT 100;" trailer on the three saved
# SQL queries in B2:B4 was un-commented to an active "LIMIT 100;" clause, and
# the selected/visible cell in the sheet was moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("B2", "B3", "B4")) {
    $cell = $ws.Range($addr)
    $row = $cell.EntireRow
    $origHeight = $row.RowHeight

    $text = $cell.Value2
    $oldTail = "--LIMIT 100;"
    if ($text.EndsWith($oldTail)) {
        $newText = $text.Substring(0, $text.Length - $oldTail.Length) + "LIMIT 100;"
        $cell.Value2 = $newText
    }

    # Un-commenting the LIMIT clause shouldn't change the (already maxed-out)
    # row height used for these long, wrapped query cells.
    $row.RowHeight = $origHeight
}

# Move the active selection / scrolled position from D3 to C4.
[void]$ws.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 2
